$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '68.034.94'
$ws.Range('E2').Value = '  +1.26%  '
Set-TextValue $ws.Range('D3') '2.635.12'
$ws.Range('E3').Value = '  +0.40%  '
$ws.Range('E4').Value = '  +0.03%  '
Set-TextValue $ws.Range('D5') '597.76'
$ws.Range('E5').Value = '  +0.10%  '
Set-TextValue $ws.Range('D6') '154.03'
$ws.Range('E6').Value = '  +1.27%  '
$ws.Range('E7').Value = '  +0.01%  '
Set-TextValue $ws.Range('D8') '0.550'
$ws.Range('E8').Value = '  -0.85%  '
Set-TextValue $ws.Range('D9') '2.634.24'
$ws.Range('E9').Value = '  +0.38%  '
$ws.Range('E10').Value = '  +10.35%  '
$ws.Range('E11').Value = '  -0.61%  '
Set-TextValue $ws.Range('D12') '5.21'
$ws.Range('E12').Value = '  +0.84%  '
Set-TextValue $ws.Range('D13') '0.347'
$ws.Range('E13').Value = '  -0.22%  '
Set-TextValue $ws.Range('D14') '27.70'
$ws.Range('E14').Value = '  +0.75%  '
Set-TextValue $ws.Range('D15') '0.0000188'
$ws.Range('E15').Value = '  +4.16%  '
Set-TextValue $ws.Range('D16') '3.112.91'
$ws.Range('E16').Value = '  +0.28%  '
Set-TextValue $ws.Range('D17') '67.881.30'
$ws.Range('E17').Value = '  +1.04%  '
Set-TextValue $ws.Range('D18') '2.637.59'
$ws.Range('E18').Value = '  +0.38%  '
Set-TextValue $ws.Range('D19') '374.51'
$ws.Range('E19').Value = '  +3.22%  '
Set-TextValue $ws.Range('D20') '11.37'
$ws.Range('E21').Value = '  +0.07%  '
Set-TextValue $ws.Range('D22') '4.25'
$ws.Range('E22').Value = '  -0.84%  '
Set-TextValue $ws.Range('D23') '4.82'
$ws.Range('E23').Value = '  -1.31%  '
$ws.Range('E24').Value = '  -2.25%  '
Set-TextValue $ws.Range('D25') '72.39'
$ws.Range('E25').Value = '  +2.12%  '
$ws.Range('E26').Value = '  +0.20%  '
Set-TextValue $ws.Range('D27') '9.98'
$ws.Range('E28').Value = '  +2.94%  '
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('E30').Value = '  -3.17%  '
Set-TextValue $ws.Range('D31') '576.97'
$ws.Range('E31').Value = '  +1.23%  '
Set-TextValue $ws.Range('D32') '1.40'
$ws.Range('E32').Value = '  +0.96%  '
Set-TextValue $ws.Range('D33') '7.85'
$ws.Range('E33').Value = '  +0.69%  '
$ws.Range('E34').Value = '  +0.04%  '
Set-TextValue $ws.Range('D35') '1.00'
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('E36').Value = '  -1.45%  '
$ws.Range('E37').Value = '  +0.34%  '
Set-TextValue $ws.Range('D38') '158.29'
$ws.Range('E38').Value = '  +0.78%  '
Set-TextValue $ws.Range('D39') '19.17'
$ws.Range('E39').Value = '  +0.13%  '
Set-TextValue $ws.Range('D40') '1.91'
$ws.Range('E40').Value = '  +5.66%  '
$ws.Range('E41').Value = '  +0.55%  '
$ws.Range('E42').Value = '  +1.95%  '
Set-TextValue $ws.Range('D43') '2.63'
$ws.Range('E43').Value = '  +4.24%  '
$ws.Range('B44').Value = 'WhiteBITCoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue $ws.Range('D44') '17.14'
$ws.Range('E44').Value = '  +4.82%  '
$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Range('D45') '0.0₆0320'
$ws.Range('E45').Value = '  +13.05%  '
$ws.Range('E46').Value = '  +0.02%  '
Set-TextValue $ws.Range('D47') '40.50'
$ws.Range('E47').Value = '  -1.74%  '
Set-TextValue $ws.Range('D48') '155.56'
$ws.Range('E48').Value = '  +0.10%  '
$ws.Range('E49').Value = '  -0.46%  '
Set-TextValue $ws.Range('D50') '22.05'
$ws.Range('E50').Value = '  +8.04%  '
Set-TextValue $ws.Range('D51') '1.70'
